$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: D=4 (Celular), G=7 (filial_id), H=8 (idCidade), N=14 (idClienteIXC), O=15 (statusInsercao), P=16 (logRetorno)
# Rows 2-9 correspond to data rows in the sheet

# --- Row 2 ---
$ws.Cells.Item(2,4).Value = "'99999999999"
$ws.Cells.Item(2,7).Value = "'35"
$ws.Cells.Item(2,8).Value = "'1659"
$ws.Cells.Item(2,14).Value = "'117636"
$ws.Cells.Item(2,16).Value = '{''type'': ''success'', ''message'': ''Registro inserido com sucesso!'', ''id'': ''117636'', ''atualiza_campos'': [{''tipo'': ''r'', ''campo'': ''ativo'', ''valor'': ''S''}, {''tipo'': ''i'', ''campo'': ''data_cadastro'', ''valor'': ''26/12/2024''}, {''tipo'': ''i'', ''campo'': ''filial_id'', ''valor'': ''35''}, {''tipo'': ''i'', ''campo'': ''latitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''longitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''id_conta'', ''valor'': ''919624''}, {''tipo'': ''d'', ''campo'': ''crm_data_vencemos'', ''valor'': ''''}, {''tipo'': ''r'', ''campo'': ''convert_cliente_forn'', ''valor'': ''''}, {''tipo'': ''d'', ''campo'': ''crm_data_perdemos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_viabilidade'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_porta_disponivel'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_abortamos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_negociando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_apresentando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sondagem'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_novo'', ''valor'': ''''}]}'

# --- Row 3 ---
$ws.Cells.Item(3,4).Value = "'99999999999"
$ws.Cells.Item(3,7).Value = "'35"
$ws.Cells.Item(3,8).Value = "'1659"
$ws.Cells.Item(3,14).Value = "'117637"
$ws.Cells.Item(3,16).Value = '{''type'': ''success'', ''message'': ''Registro inserido com sucesso!'', ''id'': ''117637'', ''atualiza_campos'': [{''tipo'': ''r'', ''campo'': ''ativo'', ''valor'': ''S''}, {''tipo'': ''i'', ''campo'': ''data_cadastro'', ''valor'': ''26/12/2024''}, {''tipo'': ''i'', ''campo'': ''filial_id'', ''valor'': ''35''}, {''tipo'': ''i'', ''campo'': ''latitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''longitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''id_conta'', ''valor'': ''919625''}, {''tipo'': ''d'', ''campo'': ''crm_data_vencemos'', ''valor'': ''''}, {''tipo'': ''r'', ''campo'': ''convert_cliente_forn'', ''valor'': ''''}, {''tipo'': ''d'', ''campo'': ''crm_data_perdemos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_viabilidade'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_porta_disponivel'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_abortamos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_negociando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_apresentando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sondagem'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_novo'', ''valor'': ''''}]}'

# --- Row 4 ---
$ws.Cells.Item(4,4).Value = "'99999999999"
$ws.Cells.Item(4,7).Value = "'35"
$ws.Cells.Item(4,8).Value = "'1659"
$ws.Cells.Item(4,14).Value = "'117638"
$ws.Cells.Item(4,16).Value = '{''type'': ''success'', ''message'': ''Registro inserido com sucesso!'', ''id'': ''117638'', ''atualiza_campos'': [{''tipo'': ''r'', ''campo'': ''ativo'', ''valor'': ''S''}, {''tipo'': ''i'', ''campo'': ''data_cadastro'', ''valor'': ''26/12/2024''}, {''tipo'': ''i'', ''campo'': ''filial_id'', ''valor'': ''35''}, {''tipo'': ''i'', ''campo'': ''latitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''longitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''id_conta'', ''valor'': ''919626''}, {''tipo'': ''d'', ''campo'': ''crm_data_vencemos'', ''valor'': ''''}, {''tipo'': ''r'', ''campo'': ''convert_cliente_forn'', ''valor'': ''''}, {''tipo'': ''d'', ''campo'': ''crm_data_perdemos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_viabilidade'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_porta_disponivel'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_abortamos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_negociando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_apresentando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sondagem'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_novo'', ''valor'': ''''}]}'

# --- Row 5 ---
$ws.Cells.Item(5,4).Value = "'99999999999"
$ws.Cells.Item(5,7).Value = "'35"
$ws.Cells.Item(5,8).Value = "'1659"
$ws.Cells.Item(5,14).Value = "'117639"
$ws.Cells.Item(5,16).Value = '{''type'': ''success'', ''message'': ''Registro inserido com sucesso!'', ''id'': ''117639'', ''atualiza_campos'': [{''tipo'': ''r'', ''campo'': ''ativo'', ''valor'': ''S''}, {''tipo'': ''i'', ''campo'': ''data_cadastro'', ''valor'': ''26/12/2024''}, {''tipo'': ''i'', ''campo'': ''filial_id'', ''valor'': ''35''}, {''tipo'': ''i'', ''campo'': ''latitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''longitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''id_conta'', ''valor'': ''919627''}, {''tipo'': ''d'', ''campo'': ''crm_data_vencemos'', ''valor'': ''''}, {''tipo'': ''r'', ''campo'': ''convert_cliente_forn'', ''valor'': ''''}, {''tipo'': ''d'', ''campo'': ''crm_data_perdemos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_viabilidade'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_porta_disponivel'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_abortamos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_negociando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_apresentando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sondagem'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_novo'', ''valor'': ''''}]}'

# --- Row 6 ---
$ws.Cells.Item(6,4).Value = "'99999999999"
$ws.Cells.Item(6,7).Value = "'35"
$ws.Cells.Item(6,8).Value = "'1659"
$ws.Cells.Item(6,14).Value = "'117640"
$ws.Cells.Item(6,15).Value = "sucesso"
$ws.Cells.Item(6,16).Value = '{''type'': ''success'', ''message'': ''Registro inserido com sucesso!'', ''id'': ''117640'', ''atualiza_campos'': [{''tipo'': ''r'', ''campo'': ''ativo'', ''valor'': ''S''}, {''tipo'': ''i'', ''campo'': ''data_cadastro'', ''valor'': ''26/12/2024''}, {''tipo'': ''i'', ''campo'': ''filial_id'', ''valor'': ''35''}, {''tipo'': ''i'', ''campo'': ''latitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''longitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''id_conta'', ''valor'': ''919628''}, {''tipo'': ''d'', ''campo'': ''crm_data_vencemos'', ''valor'': ''''}, {''tipo'': ''r'', ''campo'': ''convert_cliente_forn'', ''valor'': ''''}, {''tipo'': ''d'', ''campo'': ''crm_data_perdemos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_viabilidade'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_porta_disponivel'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_abortamos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_negociando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_apresentando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sondagem'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_novo'', ''valor'': ''''}]}'

# --- Row 7 ---
$ws.Cells.Item(7,4).Value = "'99999999999"
$ws.Cells.Item(7,7).Value = "'35"
$ws.Cells.Item(7,8).Value = "'1659"
$ws.Cells.Item(7,14).Value = "'117641"
$ws.Cells.Item(7,16).Value = '{''type'': ''success'', ''message'': ''Registro inserido com sucesso!'', ''id'': ''117641'', ''atualiza_campos'': [{''tipo'': ''r'', ''campo'': ''ativo'', ''valor'': ''S''}, {''tipo'': ''i'', ''campo'': ''data_cadastro'', ''valor'': ''26/12/2024''}, {''tipo'': ''i'', ''campo'': ''filial_id'', ''valor'': ''35''}, {''tipo'': ''i'', ''campo'': ''latitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''longitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''id_conta'', ''valor'': ''919629''}, {''tipo'': ''d'', ''campo'': ''crm_data_vencemos'', ''valor'': ''''}, {''tipo'': ''r'', ''campo'': ''convert_cliente_forn'', ''valor'': ''''}, {''tipo'': ''d'', ''campo'': ''crm_data_perdemos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_viabilidade'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_porta_disponivel'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_abortamos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_negociando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_apresentando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sondagem'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_novo'', ''valor'': ''''}]}'

# --- Row 8 ---
$ws.Cells.Item(8,4).Value = "'99999999999"
$ws.Cells.Item(8,7).Value = "'35"
$ws.Cells.Item(8,8).Value = "'1659"
$ws.Cells.Item(8,14).Value = "'117642"
$ws.Cells.Item(8,15).Value = "sucesso"
$ws.Cells.Item(8,16).Value = '{''type'': ''success'', ''message'': ''Registro inserido com sucesso!'', ''id'': ''117642'', ''atualiza_campos'': [{''tipo'': ''r'', ''campo'': ''ativo'', ''valor'': ''S''}, {''tipo'': ''i'', ''campo'': ''data_cadastro'', ''valor'': ''26/12/2024''}, {''tipo'': ''i'', ''campo'': ''filial_id'', ''valor'': ''35''}, {''tipo'': ''i'', ''campo'': ''latitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''longitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''id_conta'', ''valor'': ''919630''}, {''tipo'': ''d'', ''campo'': ''crm_data_vencemos'', ''valor'': ''''}, {''tipo'': ''r'', ''campo'': ''convert_cliente_forn'', ''valor'': ''''}, {''tipo'': ''d'', ''campo'': ''crm_data_perdemos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_viabilidade'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_porta_disponivel'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_abortamos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_negociando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_apresentando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sondagem'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_novo'', ''valor'': ''''}]}'

# --- Row 9 ---
$ws.Cells.Item(9,4).Value = "'99999999999"
$ws.Cells.Item(9,7).Value = "'35"
$ws.Cells.Item(9,8).Value = "'1659"
$ws.Cells.Item(9,14).Value = "'117643"
$ws.Cells.Item(9,15).Value = "sucesso"
$ws.Cells.Item(9,16).Value = '{''type'': ''success'', ''message'': ''Registro inserido com sucesso!'', ''id'': ''117643'', ''atualiza_campos'': [{''tipo'': ''r'', ''campo'': ''ativo'', ''valor'': ''S''}, {''tipo'': ''i'', ''campo'': ''data_cadastro'', ''valor'': ''26/12/2024''}, {''tipo'': ''i'', ''campo'': ''filial_id'', ''valor'': ''35''}, {''tipo'': ''i'', ''campo'': ''latitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''longitude'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''id_conta'', ''valor'': ''919631''}, {''tipo'': ''d'', ''campo'': ''crm_data_vencemos'', ''valor'': ''''}, {''tipo'': ''r'', ''campo'': ''convert_cliente_forn'', ''valor'': ''''}, {''tipo'': ''d'', ''campo'': ''crm_data_perdemos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_viabilidade'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sem_porta_disponivel'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_abortamos'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_negociando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_apresentando'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_sondagem'', ''valor'': ''''}, {''tipo'': ''i'', ''campo'': ''crm_data_novo'', ''valor'': ''''}]}'
